$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-10 18:25:09"
}
